# Rename the AHB-Diff column headers so that the "_old"/"_new" suffixes
# are replaced with the respective input-file format version suffixes
# ("_FV2404" / "_FV2410"), then turn the sheet into a proper Excel Table
# (Table1) over A1:U92 and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldSuffixNames = @(
  "Segmentname_old",
  "Segmentgruppe_old",
  "Segment_old",
  "Datenelement_old",
  "Segment ID_old",
  "Code_old",
  "Qualifier_old",
  "Beschreibung_old",
  "Bedingungsausdruck_old",
  "Bedingung_old"
)
$fv2404Names = @(
  "Segmentname_FV2404",
  "Segmentgruppe_FV2404",
  "Segment_FV2404",
  "Datenelement_FV2404",
  "Segment ID_FV2404",
  "Code_FV2404",
  "Qualifier_FV2404",
  "Beschreibung_FV2404",
  "Bedingungsausdruck_FV2404",
  "Bedingung_FV2404"
)

$newSuffixNames = @(
  "Segmentname_new",
  "Segmentgruppe_new",
  "Segment_new",
  "Datenelement_new",
  "Segment ID_new",
  "Code_new",
  "Qualifier_new",
  "Beschreibung_new",
  "Bedingungsausdruck_new",
  "Bedingung_new"
)
$fv2410Names = @(
  "Segmentname_FV2410",
  "Segmentgruppe_FV2410",
  "Segment_FV2410",
  "Datenelement_FV2410",
  "Segment ID_FV2410",
  "Code_FV2410",
  "Qualifier_FV2410",
  "Beschreibung_FV2410",
  "Bedingungsausdruck_FV2410",
  "Bedingung_FV2410"
)

# Columns A-J (1-10): "_old" -> "_FV2404"
for ($i = 0; $i -lt $oldSuffixNames.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $fv2404Names[$i]
}
# Column K (11) is "diff" and stays as-is.
# Columns L-U (12-21): "_new" -> "_FV2410"
for ($i = 0; $i -lt $newSuffixNames.Length; $i++) {
  $ws.Cells.Item(1, $i + 12).Value = $fv2410Names[$i]
}

$hdr = $ws.Range("A1:U1")

# Stash the header row's existing formatting on a scratch row far below the
# used range so we can restore it byte-for-byte after the table is created
# (creating a ListObject on top of an already-formatted header otherwise
# makes Excel freeze that formatting into a brand-new header-row dxf, which
# the target workbook does not have).
$template = $ws.Range("A200:U200")
$hdr.Copy()
$template.PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$hdr.ClearFormats()

# Turn the used range into a proper Excel Table.
$rng = $ws.Range("A1:U92")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.TableStyle = ""

# Restore the original header formatting (bold, fill, border, centered,
# wrapped) from the scratch copy, then discard the scratch row.
$template.Copy()
$hdr.PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Rows.Item(200).Delete()

# Freeze the header row.
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

Write-Output "done"
